# Insert two new weekly data rows ("Fruta / hortaliza, semanal") right
# before the current row 1107, shifting the existing rows 1107-1234 down
# to 1109-1236. The two new rows share the same constant columns
# (Mercado/Region/Categoria/etc.) as the rest of the block and only carry
# new Fecha / price values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert two blank rows at 1107/1108; everything below (old
# 1107..1234) slides down to 1109..1236.
$ws.Range("A1107:A1108").EntireRow.Insert()

# The row that used to be 1107 is now at 1109 - reuse its constant columns
# (A,B,C,E,F,G,H,N,O,Q,R) as the template for both freshly inserted rows.
$tmplRow = 1109

# ---- New row 1107: "Primera" ----
$ws.Cells.Item(1107,1).Value  = $ws.Cells.Item($tmplRow,1).Value()
$ws.Cells.Item(1107,2).Value  = $ws.Cells.Item($tmplRow,2).Value()
$ws.Cells.Item(1107,3).Value  = $ws.Cells.Item($tmplRow,3).Value()
$ws.Cells.Item(1107,4).Value  = 45212
$ws.Cells.Item(1107,5).Value  = $ws.Cells.Item($tmplRow,5).Value()
$ws.Cells.Item(1107,6).Value  = $ws.Cells.Item($tmplRow,6).Value()
$ws.Cells.Item(1107,7).Value  = $ws.Cells.Item($tmplRow,7).Value()
$ws.Cells.Item(1107,8).Value  = $ws.Cells.Item($tmplRow,8).Value()
$ws.Cells.Item(1107,9).Value  = "Primera"
$ws.Cells.Item(1107,10).Value = 2160
$ws.Cells.Item(1107,11).Value = 700
$ws.Cells.Item(1107,12).Value = 800
$ws.Cells.Item(1107,13).Value = 750
$ws.Cells.Item(1107,14).Value = $ws.Cells.Item($tmplRow,14).Value()
$ws.Cells.Item(1107,15).Value = $ws.Cells.Item($tmplRow,15).Value()
$ws.Cells.Item(1107,16).Value = 750
$ws.Cells.Item(1107,17).Value = $ws.Cells.Item($tmplRow,17).Value()
$ws.Cells.Item(1107,18).Value = $ws.Cells.Item($tmplRow,18).Value()

# ---- New row 1108: "Segunda" ----
$ws.Cells.Item(1108,1).Value  = $ws.Cells.Item($tmplRow,1).Value()
$ws.Cells.Item(1108,2).Value  = $ws.Cells.Item($tmplRow,2).Value()
$ws.Cells.Item(1108,3).Value  = $ws.Cells.Item($tmplRow,3).Value()
$ws.Cells.Item(1108,4).Value  = 45212
$ws.Cells.Item(1108,5).Value  = $ws.Cells.Item($tmplRow,5).Value()
$ws.Cells.Item(1108,6).Value  = $ws.Cells.Item($tmplRow,6).Value()
$ws.Cells.Item(1108,7).Value  = $ws.Cells.Item($tmplRow,7).Value()
$ws.Cells.Item(1108,8).Value  = $ws.Cells.Item($tmplRow,8).Value()
$ws.Cells.Item(1108,9).Value  = "Segunda"
$ws.Cells.Item(1108,10).Value = 900
$ws.Cells.Item(1108,11).Value = 600
$ws.Cells.Item(1108,12).Value = 600
$ws.Cells.Item(1108,13).Value = 600
$ws.Cells.Item(1108,14).Value = $ws.Cells.Item($tmplRow,14).Value()
$ws.Cells.Item(1108,15).Value = $ws.Cells.Item($tmplRow,15).Value()
$ws.Cells.Item(1108,16).Value = 600
$ws.Cells.Item(1108,17).Value = $ws.Cells.Item($tmplRow,17).Value()
$ws.Cells.Item(1108,18).Value = $ws.Cells.Item($tmplRow,18).Value()
